# updates to support new API from ApprovalBundle 2.2 and up
#
# The example data in this workbook is "rolling" - TimeTracking!B2:B14 are all
# formulas relative to StartDayConfig!A2 (the configured start day). Bumping
# the single start-day cell by exactly one year (365/366 days) refreshes every
# date shown in the example data set, which is what this commit does.

$wb = $excel.ActiveWorkbook

$startSheet = $wb.Worksheets.Item("StartDayConfig")
$trackingSheet = $wb.Worksheets.Item("TimeTracking")

# Move the example data's start date forward by exactly one year:
# 2024-01-01 -> 2025-01-01 (serial 45292 -> 45658). All of TimeTracking's date
# cells are formulas ("=StartDayConfig!A2[+n]") so they recompute automatically.
$startSheet.Range("A2").Value = 45658

# The workbook now opens on the TimeTracking tab instead of StartDayConfig.
$trackingSheet.Activate()
